# 21.4.3 Monalco Value Drive Tree Jeff Kwan - Redo.pptx
# Applies the geometry / text tweaks captured in the target diff:
#   - two "(unit)" labels -> "($)"  (with matching reflow/position of their text boxes)
#   - four small group / connector repositions (cost driver callouts)
#   - two auto-routed connector re-shapes
#   - one manually re-sized connector ("Connector: Elbow 2")

$EMU_PER_PT = 12700

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Group "Google Shape;66;p1" (id=66) - the "(unit)" label inside becomes "($)"
# and is nudged/enlarged; leave the group's own xfrm alone so the other
# (untouched) children keep their on-slide size.
# ---------------------------------------------------------------------------
$g66 = Get-ShapeById $s.Shapes 66
$lbl69 = Get-ShapeById $g66.GroupItems 69
$lbl69.TextFrame.TextRange.Text = "($)"
$lbl69.Left = (5863893 / $EMU_PER_PT)
$lbl69.Top = (1080597 / $EMU_PER_PT)
$lbl69.Width = (329898 / $EMU_PER_PT)
$lbl69.Height = (128377 / $EMU_PER_PT)

# ---------------------------------------------------------------------------
# Group "Google Shape;70;p1" (id=70) - same treatment for its "(unit)" label.
# ---------------------------------------------------------------------------
$g70 = Get-ShapeById $s.Shapes 70
$lbl73 = Get-ShapeById $g70.GroupItems 73
$lbl73.TextFrame.TextRange.Text = "($)"
$lbl73.Left = (5883212 / $EMU_PER_PT)
$lbl73.Top = (1068802 / $EMU_PER_PT)
$lbl73.Width = (329898 / $EMU_PER_PT)
$lbl73.Height = (128377 / $EMU_PER_PT)

# ---------------------------------------------------------------------------
# Groups 85 / 90 - simple horizontal repositioning.
# ---------------------------------------------------------------------------
$g85 = Get-ShapeById $s.Shapes 85
$g85.Left = (5659325 / $EMU_PER_PT)

$g90 = Get-ShapeById $s.Shapes 90
$g90.Left = (5668367 / $EMU_PER_PT)

# ---------------------------------------------------------------------------
# Connectors 95 / 96 - auto-routed bent connectors whose width shrinks to
# follow the groups above (their glued endpoints moved).
# ---------------------------------------------------------------------------
$conn95 = Get-ShapeById $s.Shapes 95
$conn95.Width = (618895 / $EMU_PER_PT)

$conn96 = Get-ShapeById $s.Shapes 96
$conn96.Width = (627937 / $EMU_PER_PT)

# ---------------------------------------------------------------------------
# Group "Google Shape;81;p1" (id=157) - moved and widened.
# ---------------------------------------------------------------------------
$g157 = Get-ShapeById $s.Shapes 157
$g157.Left = (5932560 / $EMU_PER_PT)
$g157.Top = (5267985 / $EMU_PER_PT)
$g157.Width = (1666403 / $EMU_PER_PT)

# ---------------------------------------------------------------------------
# "Connector: Elbow 2" (id=3) - manually re-routed/resized.
# ---------------------------------------------------------------------------
$conn3 = Get-ShapeById $s.Shapes 3
$conn3.Width = (1198383 / $EMU_PER_PT)
$conn3.Height = (322333 / $EMU_PER_PT)
